# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - update row 3 (R) with new simulation totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 399
$wsOff.Range("C3").Value = 282
$wsOff.Range("D3").Value = 110
$wsOff.Range("E3").Value = 56

# DEF sheet - update row 3 (R) with new simulation totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 448
$wsDef.Range("C3").Value = 332
$wsDef.Range("D3").Value = 100
$wsDef.Range("E3").Value = 50
$wsDef.Range("F3").Value = 8
